$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Gip"
$ws.Cells.Item(2, 3).Value = "Gipr"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.4870623333333333
$ws.Cells.Item(2, 8).Value = 1.461187
$ws.Cells.Item(2, 9).Value = 0.8073404988294784
$ws.Cells.Item(2, 10).Value = 0.8073404988294784
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.01627966666666667
$ws.Cells.Item(2, 14).Value = 0.048839
$ws.Cells.Item(2, 15).Value = 0.08178860127575414
$ws.Cells.Item(2, 16).Value = 0.08178860127575413
$ws.Cells.Item(2, 17).Value = 0.007929212432555555
$ws.Cells.Item(2, 18).Value = 0.071362911893
$ws.Cells.Item(2, 19).Value = 0.06603125015253267
$ws.Cells.Item(2, 20).Value = 0.06603125015253265

# Row 3
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Gip"
$ws.Cells.Item(3, 3).Value = "Gipr"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.4870623333333333
$ws.Cells.Item(3, 8).Value = 1.461187
$ws.Cells.Item(3, 9).Value = 0.8073404988294784
$ws.Cells.Item(3, 10).Value = 0.8073404988294784
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.1606193333333333
$ws.Cells.Item(3, 14).Value = 0.481858
$ws.Cells.Item(3, 15).Value = 0.8069471494816097
$ws.Cells.Item(3, 16).Value = 0.8069471494816097
$ws.Cells.Item(3, 17).Value = 0.07823162727177778
$ws.Cells.Item(3, 18).Value = 0.704084645446
$ws.Cells.Item(3, 19).Value = 0.6514811141915086
$ws.Cells.Item(3, 20).Value = 0.6514811141915086

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Gip"
$ws.Cells.Item(4, 3).Value = "Gipr"
$ws.Cells.Item(4, 4).Value = "Resolving-Mac"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.4870623333333333
$ws.Cells.Item(4, 8).Value = 1.461187
$ws.Cells.Item(4, 9).Value = 0.8073404988294784
$ws.Cells.Item(4, 10).Value = 0.8073404988294784
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.02214666666666667
$ws.Cells.Item(4, 14).Value = 0.06644
$ws.Cells.Item(4, 15).Value = 0.1112642492426361
$ws.Cells.Item(4, 16).Value = 0.1112642492426361
$ws.Cells.Item(4, 17).Value = 0.01078680714222222
$ws.Cells.Item(4, 18).Value = 0.09708126428
$ws.Cells.Item(4, 19).Value = 0.08982813448543725
$ws.Cells.Item(4, 20).Value = 0.08982813448543725

# Row 5
$ws.Cells.Item(5, 1).Value = "Resolving-Mac"
$ws.Cells.Item(5, 2).Value = "Gip"
$ws.Cells.Item(5, 3).Value = "Gipr"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.11623
$ws.Cells.Item(5, 8).Value = 0.34869
$ws.Cells.Item(5, 9).Value = 0.1926595011705215
$ws.Cells.Item(5, 10).Value = 0.1926595011705215
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.01627966666666667
$ws.Cells.Item(5, 14).Value = 0.048839
$ws.Cells.Item(5, 15).Value = 0.08178860127575414
$ws.Cells.Item(5, 16).Value = 0.08178860127575413
$ws.Cells.Item(5, 17).Value = 0.001892185656666667
$ws.Cells.Item(5, 18).Value = 0.01702967091
$ws.Cells.Item(5, 19).Value = 0.01575735112322147
$ws.Cells.Item(5, 20).Value = 0.01575735112322147

# Row 6
$ws.Cells.Item(6, 1).Value = "Resolving-Mac"
$ws.Cells.Item(6, 2).Value = "Gip"
$ws.Cells.Item(6, 3).Value = "Gipr"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.11623
$ws.Cells.Item(6, 8).Value = 0.34869
$ws.Cells.Item(6, 9).Value = 0.1926595011705215
$ws.Cells.Item(6, 10).Value = 0.1926595011705215
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.1606193333333333
$ws.Cells.Item(6, 14).Value = 0.481858
$ws.Cells.Item(6, 15).Value = 0.8069471494816097
$ws.Cells.Item(6, 16).Value = 0.8069471494816097
$ws.Cells.Item(6, 17).Value = 0.01866878511333333
$ws.Cells.Item(6, 18).Value = 0.16801906602
$ws.Cells.Item(6, 19).Value = 0.1554660352901012
$ws.Cells.Item(6, 20).Value = 0.1554660352901012

# Row 7
$ws.Cells.Item(7, 1).Value = "Resolving-Mac"
$ws.Cells.Item(7, 2).Value = "Gip"
$ws.Cells.Item(7, 3).Value = "Gipr"
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.11623
$ws.Cells.Item(7, 8).Value = 0.34869
$ws.Cells.Item(7, 9).Value = 0.1926595011705215
$ws.Cells.Item(7, 10).Value = 0.1926595011705215
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.02214666666666667
$ws.Cells.Item(7, 14).Value = 0.06644
$ws.Cells.Item(7, 15).Value = 0.1112642492426361
$ws.Cells.Item(7, 16).Value = 0.1112642492426361
$ws.Cells.Item(7, 17).Value = 0.002574107066666667
$ws.Cells.Item(7, 18).Value = 0.0231669636
$ws.Cells.Item(7, 19).Value = 0.02143611475719885
$ws.Cells.Item(7, 20).Value = 0.02143611475719885
